$wb = $excel.ActiveWorkbook

# --- Sheet "All Published Values": append new row 8 with the latest published rate ---
$ws = $wb.Worksheets.Item("All Published Values")

# Force the new row to be written as plain text (matches the existing rows,
# which store every value - including dates/numbers - as text) instead of
# letting Excel auto-detect dates/numbers.
$ws.Range("A8:J8").NumberFormat = "@"

$ws.Range("A8").Value = "2026-01-02"
$ws.Range("B8").Value = "2026-01-02 19:33:08"
$ws.Range("C8").Value = "697.85"
$ws.Range("D8").Value = "697.85"
$ws.Range("E8").Value = "700.79"
$ws.Range("F8").Value = "700.79"
$ws.Range("G8").Value = "702.88"
$ws.Range("H8").Value = "2026/01/02 19:33:08"
$ws.Range("I8").Value = "2026-01-02 11:38:01"
$ws.Range("J8").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Re-apply the AutoFilter so its range grows from A1:J7 to A1:J8.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
[void]$ws.Range("A1:J8").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$8"
    }
}

# --- Sheet "Daily Summary": bump the publishes count for 2026-01-02 ---
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("B4").Value = 7
